# The "Key" sheet lists field-name / ESRI-field-name / description rows.
# Three rows documenting fields that no longer exist in the data
# (IHO_SEA, TERRITORY, COUNTRY) are removed, and everything below shifts up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 = COUNTRY (remove first so row numbers for the block above stay valid)
$ws.Rows("16:16").Delete()
# Rows 12:13 = IHO_SEA, TERRITORY
$ws.Rows("12:13").Delete()

# Update the view: zoom to 120% and move the active selection to B14
$excel.ActiveWindow.Zoom = 120
$ws.Range("B14").Select()
